$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.427.79'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.674.13'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("E4").Value = '  +0.80%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5355'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2675'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06411'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07861'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.549'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.659.61'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.904.50'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5666'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8217'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.470.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.738'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '199.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.92%  '
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("E24").Value = '  +0.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1233'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.271'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.505'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05904'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.287'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.587'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.320'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.624'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9710'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.04%  '
$ws.Range("E36").Value = '  +1.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.440'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5835'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01619'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.081.42'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.925'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8674'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '104.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.813.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.09%  '
$ws.Range("E47").Value = '  -3.74%  '
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4416'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.047'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05171'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.39%  '
